# Updates the Saudi Professional League 2023-2024 odds sheet:
#  - fixes two matches per round that had been paired to the wrong fixture
#    (rows 18/19, 81/82, 92/93, 105/106 had their match-specific columns
#    swapped between the two rows sharing the same kickoff date/time)
#  - appends two newly scraped fixtures (rows 114 and 115)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowColumns($row1, $row2, $cols) {
    foreach ($col in $cols) {
        $r1 = $ws.Range($col + $row1)
        $r2 = $ws.Range($col + $row2)
        $v1 = $r1.Value()
        $v2 = $r2.Value()
        $r1.Value = $v2
        $r2.Value = $v1
    }
}

# Columns that carry match-specific data (team names, scores, odds, odds
# timestamps, match url). A-E (index/pais/torneio/temporada/data_partida)
# and K/O/S (opening odds timestamps) are identical between each paired
# row, so they are intentionally left untouched.
$matchCols = @("F", "G", "H", "I", "J", "L", "M", "N", "P", "Q", "R", "T", "U", "V")

Swap-RowColumns 18 19 $matchCols
Swap-RowColumns 81 82 $matchCols
Swap-RowColumns 92 93 $matchCols
Swap-RowColumns 105 106 $matchCols

# Append the two new fixtures as rows 114 and 115, copying the formatting
# (bold/bordered index column, datetime number format on data_partida)
# from the last existing row (113).
$ws.Range("A113:V113").Copy()
$ws.Range("A114:V115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$row114 = @{
    A = 113
    B = "saudi-arabia"
    C = "saudi-professional-league"
    D = "2023-2024"
    E = 45240.66666666666
    F = "Al Hilal"
    G = 2
    H = "Al Taawon"
    I = 0
    J = 1.23
    K = "05/11/2023 16:12"
    L = 1.43
    M = "10/11/2023 15:59"
    N = 6.98
    O = "05/11/2023 16:12"
    P = 5.54
    Q = "10/11/2023 15:59"
    R = 10.09
    S = "05/11/2023 16:12"
    T = 5.92
    U = "10/11/2023 15:59"
    V = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-hilal-al-taawon/pGQuNnci/"
}

$row115 = @{
    A = 114
    B = "saudi-arabia"
    C = "saudi-professional-league"
    D = "2023-2024"
    E = 45240.79166666666
    F = "Al Ittihad"
    G = 4
    H = "Abha"
    I = 2
    J = 1.11
    K = "04/11/2023 17:13"
    L = 1.19
    M = "10/11/2023 18:54"
    N = 10.79
    O = "04/11/2023 17:13"
    P = 7.5
    Q = "10/11/2023 18:58"
    R = 18.91
    S = "04/11/2023 17:13"
    T = 12.7
    U = "10/11/2023 18:58"
    V = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ittihad-abha/0lMLHiLB/"
}

foreach ($col in $row114.Keys) {
    $ws.Range($col + "114").Value = $row114[$col]
}

foreach ($col in $row115.Keys) {
    $ws.Range($col + "115").Value = $row115[$col]
}

Write-Host "Applied swaps and appended rows 114-115"
